$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 4): Bug/Error, Solution, Date
$ws.Range("B4").Value = "More Info button through JS  not clearing old for and displaying new"
$ws.Range("C4").Value = "CSS/HTML work around "

# Copy the date formatting from D3 so D4 reuses the existing date style
# instead of Excel auto-generating a brand-new number format.
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 43486

# Leave the selection where the author ended up after entering the row.
[void]$ws.Range("D5").Select()
